$d = $word.ActiveDocument

function Add-ParagraphBeforeTrailing($text) {
    $trailingPara = $d.Paragraphs.Last
    $r = $trailingPara.Range
    $r.Collapse(1)
    $null = $r.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $nr = $newPara.Range
    $nr.Collapse(1)
    $null = $nr.InsertBefore($text)
}

# Position at the very end of the document (right after "esfriar.")
$last = $d.Paragraphs.Last
$rng = $last.Range
$rng.Collapse(0)
$null = $rng.InsertParagraphAfter()

# New paragraph: " " + "DICAS DE PREPARO:" as two separate runs, inserted via
# raw OOXML so the leading space keeps its own run (matching the authored edit).
$newPara = $d.Paragraphs.Last
$paraRange = $newPara.Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>DICAS DE PREPARO:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $paraRange.InsertXML($xml)

# Remaining tip paragraphs, inserted right before the trailing empty paragraph
# that InsertXML left behind (so that empty paragraph ends up last, as in the diff).
Add-ParagraphBeforeTrailing("- Para um bolo mais fofo, peneire a farinha de trigo.")
Add-ParagraphBeforeTrailing("- Você pode adicionar nozes picadas à massa para um toque especial.")
